# Auto-generated edit script applying the diff changes
# Updates "想去人数" (want-to-go count) / "最低票价" (min price) counters
# and refreshes several event listing rows (name/venue/time/link/cover)
# across the 展览, 演出 and 全部类型 sheets to match the newly scraped data.
$wb = $excel.ActiveWorkbook

# ----- Sheet: 展览 -----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 10275
$ws.Range("F5").Value = 754
$ws.Range("C6").Value = "北京·【五一艺术展】奇点艺术节.ARTPHILE®2024"
$ws.Range("D6").Value = "朝阳门外大街10号 THE BOX 朝外B座"
$ws.Range("E6").Value = "2024.05.01 11:00-05.04 20:00"
$ws.Range("F6").Value = 13
$ws.Range("G6").Value = 80
$ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=84864"
$ws.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/1xX6GkXr1713958631868.jpeg"
$ws.Range("C7").Value = "北京·卡淘嘉年华·第三届球星卡交流会"
$ws.Range("D7").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Range("E7").Value = "2024.05.01 09:30-05.03 17:00"
$ws.Range("F7").Value = 206
$ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=82072"
$ws.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202402/XOTabMFt1708929919204.jpeg"
$ws.Range("C8").Value = "北京·原神x穹铁北京同人嘉年华7th"
$ws.Range("D8").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E8").Value = "2024.05.01 09:00-05.04 17:00"
$ws.Range("F8").Value = 459
$ws.Range("G8").Value = 95
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84114"
$ws.Range("I8").Value = "//i2.hdslb.com/bfs/openplatform/202404/55ApL1HY1712813894389.jpeg"
$ws.Range("C9").Value = "北京·嘉品次元派对-免费展会"
$ws.Range("D9").Value = "东坝中路38号 北京金隅嘉品Mall中庭"
$ws.Range("E9").Value = "2024.05.01 14:00-05.05 20:30"
$ws.Range("F9").Value = 431
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84171"
$ws.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/utfRydGx1713238690950.jpeg"
$ws.Range("C10").Value = "北京·国乙同好嘉年华7th"
$ws.Range("D10").Value = "北京国家会议中心 北京国家会议中心"
$ws.Range("E10").Value = "2024.05.01 09:00-05.04 17:00"
$ws.Range("F10").Value = 477
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=82391"
$ws.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202403/BGYIf9qe1709696198696.jpeg"
$ws.Range("C11").Value = "北京·排球少年同好嘉年华2nd"
$ws.Range("D11").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E11").Value = "2024.05.01 09:30-05.04 17:00"
$ws.Range("F11").Value = 269
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=84070"
$ws.Range("I11").Value = "//i1.hdslb.com/bfs/openplatform/202404/UHohfLBe1714358642039.jpeg"
$ws.Range("B12").Value = "'2024-05-01"
$ws.Range("C12").Value = "北京·第16届IJOY漫展XCGF游戏节"
$ws.Range("D12").Value = "北京国家会议中心 北京国家会议中心"
$ws.Range("E12").Value = "2024.05.01 09:00-05.04 17:00"
$ws.Range("F12").Value = 12839
$ws.Range("G12").Value = 95
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=81183"
$ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202402/H86O2Jvv1707017473134.jpeg"
$ws.Range("B14").Value = "'2024-05-03"
$ws.Range("C14").Value = "北京·知名演员 川久保拓司 专场活动"
$ws.Range("D14").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Range("E14").Value = "2024.05.03 10:30-05.03 15:00"
$ws.Range("F14").Value = 161
$ws.Range("G14").Value = 528
$ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=82897"
$ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202403/rxrJuuvX1710409029498.jpeg"
$ws.Range("B15").Value = "'2024-05-04"
$ws.Range("C15").Value = "北京·XW咒术回战only"
$ws.Range("D15").Value = "北花园路1号 超级蜂巢"
$ws.Range("E15").Value = "2024.05.04 10:00-05.04 17:00"
$ws.Range("F15").Value = 277
$ws.Range("G15").Value = 60
$ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=83570"
$ws.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202403/G9X2HmU11711703284044.jpeg"
$ws.Range("C16").Value = "北京·第16届IJOY漫展【文森个人专场见面会】"
$ws.Range("D16").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E16").Value = "2024.05.04 11:00-05.04 15:10"
$ws.Range("F16").Value = 49
$ws.Range("G16").Value = 238
$ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=83617"
$ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202404/E3ZY4mKN1711961443069.jpeg"
$ws.Range("B17").Value = "'2024-05-12"
$ws.Range("C17").Value = "北京·次元仙界会"
$ws.Range("D17").Value = "丽泽天地购物中心 丽泽天地购物中心"
$ws.Range("E17").Value = "2024.05.12 10:00-05.13 02:00"
$ws.Range("F17").Value = 194
$ws.Range("G17").Value = 49
$ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=83690"
$ws.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202404/nw4FcqlD1712029131170.png"
$ws.Range("B18").Value = "'2024-05-18"
$ws.Range("C18").Value = "北京·ICOS SP漫展04X五只猫动漫节"
$ws.Range("D18").Value = "北京电影学院影视文化产业创新园平房园区 北京五只猫娱乐Mall"
$ws.Range("E18").Value = "2024.05.18 09:00-05.19 17:00"
$ws.Range("F18").Value = 147
$ws.Range("G18").Value = 80
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=83122"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202403/3N8tBAKl1710831573887.jpeg"
$ws.Range("C19").Value = "北京·YIYOU二次元大聚会"
$ws.Range("D19").Value = "京开高速入口与京开高速交叉口西180米 北京双马文体创业园"
$ws.Range("E19").Value = "2024.05.18 10:00-05.18 18:00"
$ws.Range("F19").Value = 184
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=83129"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZhTtVA3A1710812150528.png"
$ws.Range("C20").Value = "北京·原神only3.0"
$ws.Range("D20").Value = "北花园路1号 超级蜂巢"
$ws.Range("E20").Value = "2024.05.18 10:00-05.19 17:00"
$ws.Range("F20").Value = 2769
$ws.Range("G20").Value = 68
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=81766"
$ws.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202402/Lfxwe5PO1707120983684.jpeg"
$ws.Range("C21").Value = "北京·原神only3.0——32D小神奈签售会"
$ws.Range("E21").Value = "2024.05.18 10:00-05.18 17:00"
$ws.Range("F21").Value = 48
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=82147"
$ws.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202402/lQoExxJd1709100610683.jpeg"
$ws.Range("F24").Value = 111
$ws.Range("F29").Value = 1119
$ws.Range("F30").Value = 4283
$ws.Range("F32").Value = 3838
$ws.Range("F33").Value = 891
$ws.Range("F34").Value = 2658
$ws.Range("F35").Value = 3091
$ws.Range("F36").Value = 84
$ws.Range("F37").Value = 1373
$ws.Range("F39").Value = 791
$ws.Range("F41").Value = 148
$ws.Range("F42").Value = 536
$ws.Range("F43").Value = 750
$ws.Range("G43").Value = 85
$ws.Range("F45").Value = 169
$ws.Range("F46").Value = 315
$ws.Range("G46").Value = 85
$ws.Range("F47").Value = 124
$ws.Range("F48").Value = 172
$ws.Range("F49").Value = 188
$ws.Range("G49").Value = 85

# ----- Sheet: 演出 -----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 68
$ws.Range("F8").Value = 38

# ----- Sheet: 全部类型 -----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 10275
$ws.Range("F5").Value = 754
$ws.Range("F7").Value = 459
$ws.Range("F8").Value = 431
$ws.Range("F11").Value = 12839
$ws.Range("F14").Value = 49
$ws.Range("F15").Value = 68
$ws.Range("F16").Value = 38
$ws.Range("G17").Value = 49
$ws.Range("F19").Value = 184
$ws.Range("F20").Value = 2769
$ws.Range("F22").Value = 111
$ws.Range("F27").Value = 1119
$ws.Range("F29").Value = 4283
$ws.Range("F30").Value = 3838
$ws.Range("F31").Value = 891
$ws.Range("F32").Value = 2658
$ws.Range("F33").Value = 3091
$ws.Range("F34").Value = 84
$ws.Range("F36").Value = 1373
$ws.Range("F38").Value = 791
$ws.Range("F40").Value = 148
$ws.Range("F41").Value = 536
$ws.Range("F43").Value = 750
$ws.Range("G43").Value = 85
$ws.Range("F45").Value = 169
$ws.Range("F46").Value = 315
$ws.Range("G46").Value = 85
$ws.Range("F47").Value = 124
$ws.Range("F48").Value = 172
$ws.Range("F49").Value = 188
$ws.Range("G49").Value = 85
